# Update raw data collection (sanity check, iotamotion device) plus unit correction on graph
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "MetalShape"
$ws.Range("B1").Value = "MetalNumber"
$ws.Range("C1").Value = "Validity_Check_Value_of_Translation_Error"
$ws.Range("D1").Value = "Validity_Check_Value_of_Rotation_Error"

$ws.Range("A2").Value = "hollow"
$ws.Range("B2").Value = "LC Steel"
$ws.Range("C2").Value = 0.000000993902614153660331845742
$ws.Range("D2").Value = 0.000008229527849090802202390409

$ws.Range("A3").Value = "hollow"
$ws.Range("B3").Value = "416 SS"
$ws.Range("C3").Value = 0.000000174790407552334554002567
$ws.Range("D3").Value = 0.000004352947695910384422125945

$ws.Range("A4").Value = "hollow"
$ws.Range("B4").Value = "304 SS"
$ws.Range("C4").Value = 0.000000466743207585291768445669
$ws.Range("D4").Value = 0.000003808521576142061360677367

$ws.Range("A5").Value = "hollow"
$ws.Range("B5").Value = "6061 Al"
$ws.Range("C5").Value = 0.000000787091704337074673185739
$ws.Range("D5").Value = 0.000009568988608052408003919014

$ws.Range("A6").Value = "hollow"
$ws.Range("B6").Value = "Ti Grade 5"
$ws.Range("C6").Value = 0.000000615450567995568129495327
$ws.Range("D6").Value = 0.000009279920890481894034016418

$ws.Range("A7").Value = "hollow"
$ws.Range("B7").Value = "Copper"
$ws.Range("C7").Value = 0.000000651232144250815492457431
$ws.Range("D7").Value = 0.000004637346174973210834906797

$ws.Range("A8").Value = "sheet"
$ws.Range("B8").Value = "LC Steel"
$ws.Range("C8").Value = 0.000000167609090603334432523374
$ws.Range("D8").Value = 0.000001460736534309842550772738

$ws.Range("A9").Value = "sheet"
$ws.Range("B9").Value = "304 SS"
$ws.Range("C9").Value = 0.000000463138481554189727589627
$ws.Range("D9").Value = 0.000005906038866372831749223041

$ws.Range("A10").Value = "sheet"
$ws.Range("B10").Value = "6061 Al"
$ws.Range("C10").Value = 0.000000432812374265341651705325
$ws.Range("D10").Value = 0.000004573692670541227609057756

$ws.Range("A11").Value = "sheet"
$ws.Range("B11").Value = "Copper"
$ws.Range("C11").Value = 0.000000780479244563513836276689
$ws.Range("D11").Value = 0.000002646335333995063524795827

$ws.Range("A12").Value = "solid"
$ws.Range("B12").Value = "LC Steel"
$ws.Range("C12").Value = 0.000000168079776539078728751170
$ws.Range("D12").Value = 0.000001641981700463392682790855

$ws.Range("A13").Value = "solid"
$ws.Range("B13").Value = "416 SS"
$ws.Range("C13").Value = 0.000000849761912973042658403706
$ws.Range("D13").Value = 0.000003466866103945684935127560

$ws.Range("A14").Value = "solid"
$ws.Range("B14").Value = "304 SS"
$ws.Range("C14").Value = 0.000000408640735918643763973577
$ws.Range("D14").Value = 0.000004344463691038418213347886

$ws.Range("A15").Value = "solid"
$ws.Range("B15").Value = "6061 Al"
$ws.Range("C15").Value = 0.000000377406915936643707937433
$ws.Range("D15").Value = 0.000008365137507033221019947869

$ws.Range("A16").Value = "solid"
$ws.Range("B16").Value = "Ti Grade 5"
$ws.Range("C16").Value = 0.000000584802776312902659769334
$ws.Range("D16").Value = 0.000003390276512204230967192356

$ws.Range("A17").Value = "solid"
$ws.Range("B17").Value = "Copper"
$ws.Range("C17").Value = 0.000000151291402113363108161828
$ws.Range("D17").Value = 0.000002417612091312938303253044
